$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 148-149, pushing the existing rows 148-160 down to 150-162.
$ws.Range("A148:R149").Insert()

# Populate the newly inserted row 148 (Magnum) with the new week's data.
$ws.Range("A148").Value = 2
$ws.Range("B148").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44615
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 100112031
$ws.Range("G148").Value = "Poroto verde"
$ws.Range("H148").Value = "Magnum"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 800
$ws.Range("K148").Value = 20000
$ws.Range("L148").Value = 21000
$ws.Range("M148").Value = 20500
$ws.Range("N148").Value = '$/malla 25 kilos'
$ws.Range("O148").Value = "Provincia de Limarí"
$ws.Range("P148").Value = 820
$ws.Range("Q148").Value = 25
$ws.Range("R148").Value = "Hortaliza"

# Populate the newly inserted row 149 (Sin especificar) with the new week's data.
$ws.Range("A149").Value = 2
$ws.Range("B149").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44615
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100112031
$ws.Range("G149").Value = "Poroto verde"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 500
$ws.Range("K149").Value = 23000
$ws.Range("L149").Value = 25000
$ws.Range("M149").Value = 24000
$ws.Range("N149").Value = '$/malla 25 kilos'
$ws.Range("O149").Value = "Provincia de Limarí"
$ws.Range("P149").Value = 960
$ws.Range("Q149").Value = 25
$ws.Range("R149").Value = "Hortaliza"
